$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.94"
$ws.Range("E2").Value = "'0.19%"
$ws.Range("D3").Value = "'32.25"
$ws.Range("E3").Value = "'2.12%"
$ws.Range("D4").Value = "'4.952"
$ws.Range("E4").Value = "'-2.91%"
$ws.Range("D5").Value = "'0.07646"
$ws.Range("E5").Value = "'-2.09%"
$ws.Range("D6").Value = "'1.921"
$ws.Range("E6").Value = "'-17.46%"
$ws.Range("D7").Value = "'7.836"
$ws.Range("B8").Value = "'GateToken"
$ws.Range("C8").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'3.800"
$ws.Range("E8").Value = "'-0.85%"
$ws.Range("B9").Value = "'MXToken"
$ws.Range("C9").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9172"
$ws.Range("E9").Value = "'0.31%"
$ws.Range("B10").Value = "'WazirX"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1746"
$ws.Range("E10").Value = "'0.01%"
$ws.Range("B11").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.07763"
$ws.Range("E11").Value = "'3.12%"
$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08528"
$ws.Range("E12").Value = "'-6.80%"
$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03170"
$ws.Range("E13").Value = "'2.19%"
$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09985"
$ws.Range("E14").Value = "'-0.33%"
$ws.Range("B15").Value = "'BitForexToken"
$ws.Range("C15").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001518"
$ws.Range("E15").Value = "'0.80%"
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005916"
$ws.Range("E16").Value = "'1.08%"
$ws.Range("B17").Value = "'LEO"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.465"
$ws.Range("E17").Value = "'-0.42%"
$ws.Range("E18").Value = "'-4.23%"
$ws.Range("D19").Value = "'0.3350"
$ws.Range("E19").Value = "'1.82%"
$ws.Range("D20").Value = "'0.1327"
$ws.Range("E20").Value = "'-0.83%"
$ws.Range("D21").Value = "'4.290"
$ws.Range("E21").Value = "'6.60%"
$ws.Range("D23").Value = "'0.04516"
$ws.Range("E23").Value = "'-1.62%"
$ws.Range("D24").Value = "'0.001221"
$ws.Range("E24").Value = "'-2.40%"
$ws.Range("D25").Value = "'0.004386"
$ws.Range("E25").Value = "'-1.59%"
$ws.Range("E26").Value = "'0.03%"
$ws.Range("D39").Value = "'0.01696"
$ws.Range("E39").Value = "'-4.44%"
$ws.Range("D40").Value = "'0.04695"
$ws.Range("E40").Value = "'-1.79%"
$ws.Range("D41").Value = "'0.007460"
$ws.Range("E41").Value = "'0.64%"
$ws.Range("E42").Value = "'-0.58%"
$ws.Range("E43").Value = "'6.42%"
$ws.Range("D44").Value = "'0.01048"
$ws.Range("E44").Value = "'2.38%"
$ws.Range("D45").Value = "'0.00006242"
$ws.Range("E45").Value = "'0.54%"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("D47").Value = "'0.8234"
$ws.Range("E47").Value = "'10.49%"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("D50").Value = "'0.0002003"
$ws.Range("E50").Value = "'0.05%"
